$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the date column as literal text ("01/12/2023"), not an auto-converted
# date serial number, by forcing a text number format before assigning it
# (the equivalent of typing a leading apostrophe in the UI).
$ws.Range("A2:A3").NumberFormat = "@"

# Row 2: date, temperature unit ("C"), humidity unit ("%").
$ws.Range("A2").Value = "01/12/2023"
$ws.Range("J2").Value = "C"
$ws.Range("K2").Value = "%"

# Row 3: same date/unit values, plus the remaining mid-row cells (B3:I3, L3)
# are touched but left blank, matching the placeholder cells in the source
# row that carry no text.
$ws.Range("B3:I3").NumberFormat = "@"
$ws.Range("L3").NumberFormat = "@"
$ws.Range("A3").Value = "01/12/2023"
$ws.Range("J3").Value = "C"
$ws.Range("K3").Value = "%"
